$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2: PCA loadings table - header cell (row 1, col 3) gets reworded.
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$tbl = $slide2.Shapes.Item(1).Table
$headerCell = $tbl.Cell(1, 3)
$headerRange = $headerCell.Shape.TextFrame.TextRange
$headerRange.Text = "Starvation " + [char]0x2013 + "desiccation/thermeral tolerance trade off(PC2, 22.25%)"

# ---------------------------------------------------------------------------
# Slide 5: new caption textbox above the right-hand scatterplot.
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$box5 = $slide5.Shapes.AddTextbox(1, 364.9033070866142, 80.93748031496062, 598.5133858267717, 29.081259842519685)
$box5.TextFrame.TextRange.Text = "High starvation hardiness and low thermal tolerance, low desiccation hardiness"
$box5.Fill.Visible = 0
$box5.TextFrame.WordWrap = 0
$box5.TextFrame.AutoSize = 1

# ---------------------------------------------------------------------------
# Slide 6: three new caption textboxes around the big scatterplot figure.
# ---------------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)

# "Left panel: ..." caption (bottom-left, 3 runs)
$boxLeft = $slide6.Shapes.AddTextbox(1, 177.97755905511812, 475.68535433070866, 582.9679527559055, 29.081259842519685)
$trLeft = $boxLeft.TextFrame.TextRange
$trLeft.Text = "Left panel: Low growth and high cold tolerance are at the edge of the range"
$trLeft.Characters(1, 33).Text = "Left panel: Low growth and high "
$trLeft.Characters(34, 19).Text = "cold tolerance are "
$trLeft.Characters(53, 24).Text = "at the edge of the range"
$boxLeft.Fill.Visible = 0
$boxLeft.TextFrame.WordWrap = 0
$boxLeft.TextFrame.AutoSize = 1

# "Right panel: ..." caption (bottom-right, 3 runs)
$boxRight = $slide6.Shapes.AddTextbox(1, 177.97755905511812, 500.1462992125984, 646.7044094488189, 29.081259842519685)
$trRight = $boxRight.TextFrame.TextRange
$trRight.Text = "Right panel: High growth and low stress hardiness is associated with lower precitation "
$trRight.Characters(1, 76).Text = "Right panel: High growth and low stress hardiness is associated with lower "
$trRight.Characters(77, 11).Text = "precitation"
$trRight.Characters(88, 1).Text = " "
$boxRight.Fill.Visible = 0
$boxRight.TextFrame.WordWrap = 0
$boxRight.TextFrame.AutoSize = 1

# "More cold tolerant, ..." caption (top, single run)
$boxTop = $slide6.Shapes.AddTextbox(1, 8.08992125984252, 74.1896062992126, 717.9988976377953, 29.081259842519685)
$boxTop.TextFrame.TextRange.Text = "More cold tolerant, higher growth, lower heat knock down , lower starvation, lower desiccation"
$boxTop.Fill.Visible = 0
$boxTop.TextFrame.WordWrap = 0
$boxTop.TextFrame.AutoSize = 1
